$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.936384891369244
$ws.Range("C2").Value = 0.3063928068183088
$ws.Range("E2").Value = 0.06041624720933214
$ws.Range("F2").Value = 3.503324868206647
$ws.Range("G2").Value = 0.002588067574074
$ws.Range("I2").Value = 2.314698982480635
$ws.Range("J2").Value = 0.1778744747948195
$ws.Range("L2").Value = 0.298820149118491
$ws.Range("M2").Value = 0.4143508551076067
$ws.Range("B3").Value = 1.863314276690176
$ws.Range("C3").Value = 0.2735367505877662
$ws.Range("E3").Value = 0.05892210267657561
$ws.Range("F3").Value = 3.468298162815032
$ws.Range("G3").Value = 0.002593931998298218
$ws.Range("I3").Value = 2.291330996386179
$ws.Range("J3").Value = 0.1779329282893585
$ws.Range("L3").Value = 0.2989477356285306
$ws.Range("M3").Value = 0.4052213421678204
$ws.Range("B4").Value = 1.819783652279853
$ws.Range("C4").Value = 0.253504437051447
$ws.Range("E4").Value = 0.05799221195438342
$ws.Range("F4").Value = 3.44829744513612
$ws.Range("G4").Value = 0.002597720999803761
$ws.Range("I4").Value = 2.27790491901068
$ws.Range("J4").Value = 0.1779988299064499
$ws.Range("L4").Value = 0.2991545467433809
$ws.Range("M4").Value = 0.3998826666199449
$ws.Range("B5").Value = 1.802379702974804
$ws.Range("C5").Value = 0.2453757656786877
$ws.Range("E5").Value = 0.0576100608997816
$ws.Range("F5").Value = 3.440524291685222
$ws.Range("G5").Value = 0.002599312544576344
$ws.Range("I5").Value = 2.272664185406967
$ws.Range("J5").Value = 0.1780331580372625
$ws.Range("L5").Value = 0.299271248514529
$ws.Range("M5").Value = 0.3977742518243872
$ws.Range("B6").Value = 1.799510015891599
$ws.Range("C6").Value = 0.2440280730220081
$ws.Range("E6").Value = 0.05754640876958916
$ws.Range("F6").Value = 3.439256314338536
$ws.Range("G6").Value = 0.002599579692815529
$ws.Range("I6").Value = 2.271807848472406
$ws.Range("J6").Value = 0.1780393073716944
$ws.Range("L6").Value = 0.2992925886475462
$ws.Range("M6").Value = 0.3974282071665201
$ws.Range("B7").Value = 1.819547580331914
$ws.Range("C7").Value = 0.2533946717975368
$ws.Range("E7").Value = 0.05798707123815561
$ws.Range("F7").Value = 3.44819108780213
$ws.Range("G7").Value = 0.002597742271359734
$ws.Range("I7").Value = 2.277833309136398
$ws.Range("J7").Value = 0.1779992627102409
$ws.Range("L7").Value = 0.2991559891784163
$ws.Range("M7").Value = 0.3998539599175146
$ws.Range("B8").Value = 1.910912674900828
$ws.Range("C8").Value = 0.295034124947648
$ws.Range("E8").Value = 0.05990361089645369
$ws.Range("F8").Value = 3.490934294755789
$ws.Range("G8").Value = 0.002590050662542154
$ws.Range("I8").Value = 2.30644953425805
$ws.Range("J8").Value = 0.1778883562670153
$ws.Range("L8").Value = 0.2988375467958662
$ws.Range("M8").Value = 0.4111476190532244
$ws.Range("B9").Value = 2.100712076733259
$ws.Range("C9").Value = 0.3778564932026711
$ws.Range("E9").Value = 0.06356642392202438
$ws.Range("F9").Value = 3.586770262124418
$ws.Range("G9").Value = 0.002576453299105498
$ws.Range("I9").Value = 2.369945884706397
$ws.Range("J9").Value = 0.1779123483081086
$ws.Range("L9").Value = 0.2992275693183473
$ws.Range("M9").Value = 0.4354128874429364
$ws.Range("B10").Value = 2.246715398947856
$ws.Range("C10").Value = 0.4394868741865139
$ws.Range("E10").Value = 0.06620430902303198
$ws.Range("F10").Value = 3.664612951956002
$ws.Range("G10").Value = 0.002567358482274465
$ws.Range("I10").Value = 2.421191874673198
$ws.Range("J10").Value = 0.1780818734196359
$ws.Range("L10").Value = 0.3001265250826819
$ws.Range("M10").Value = 0.4545363341630306
$ws.Range("B11").Value = 2.314578300121752
$ws.Range("C11").Value = 0.4677091088549901
$ws.Range("E11").Value = 0.06739391004777673
$ws.Range("F11").Value = 3.701663233406066
$ws.Range("G11").Value = 0.002563413114322244
$ws.Range("I11").Value = 2.445524058359425
$ws.Range("J11").Value = 0.1781929955734469
$ws.Range("L11").Value = 0.3006672385520801
$ws.Range("M11").Value = 0.4635186771627389
$ws.Range("B12").Value = 2.340485208165944
$ws.Range("C12").Value = 0.4784240624601921
$ws.Range("E12").Value = 0.06784297435595832
$ws.Range("F12").Value = 3.715930672638081
$ws.Range("G12").Value = 0.002561946528245836
$ws.Range("I12").Value = 2.454886321507246
$ws.Range("J12").Value = 0.1782400472831966
$ws.Range("L12").Value = 0.3008908351543482
$ws.Range("M12").Value = 0.4669608045582549
$ws.Range("B13").Value = 2.334896395559042
$ws.Range("C13").Value = 0.4761151523899798
$ws.Range("E13").Value = 0.06774632180116313
$ws.Range("F13").Value = 3.71284734489663
$ws.Range("G13").Value = 0.002562261165812101
$ws.Range("I13").Value = 2.452863375659447
$ws.Range("J13").Value = 0.1782296913989896
$ws.Range("L13").Value = 0.3008418434185742
$ws.Range("M13").Value = 0.4662176699710088
$ws.Range("B14").Value = 2.316705489200729
$ws.Range("C14").Value = 0.4685900713606657
$ws.Range("E14").Value = 0.06743088266146913
$ws.Range("F14").Value = 3.702832259223641
$ws.Range("G14").Value = 0.00256329190841926
$ws.Range("I14").Value = 2.446291318732335
$ws.Range("J14").Value = 0.1781967664224382
$ws.Range("L14").Value = 0.3006852570114091
$ws.Range("M14").Value = 0.4638010470727423
$ws.Range("B15").Value = 2.30559023947967
$ws.Range("C15").Value = 0.4639843936666352
$ws.Range("E15").Value = 0.06723748566596655
$ws.Range("F15").Value = 3.696728680964043
$ws.Range("G15").Value = 0.002563926836469474
$ws.Range("I15").Value = 2.442285087976089
$ws.Range("J15").Value = 0.1781772489002975
$ws.Range("L15").Value = 0.300591793660189
$ws.Range("M15").Value = 0.4623260978623875
$ws.Range("B16").Value = 2.242309554653104
$ws.Range("C16").Value = 0.4376463198436227
$ws.Range("E16").Value = 0.06612636376896397
$ws.Range("F16").Value = 3.662224749425917
$ws.Range("G16").Value = 0.002567620169584825
$ws.Range("I16").Value = 2.419622364763626
$ws.Range("J16").Value = 0.1780753026082635
$ws.Range("L16").Value = 0.3000938304750136
$ws.Range("M16").Value = 0.4539550125714982
$ws.Range("B17").Value = 2.203859676072909
$ws.Range("C17").Value = 0.4215371932747303
$ws.Range("E17").Value = 0.06544212420838669
$ws.Range("F17").Value = 3.641478674262714
$ws.Range("G17").Value = 0.002569934950461477
$ws.Range("I17").Value = 2.40598186051389
$ws.Range("J17").Value = 0.1780215296199366
$ws.Range("L17").Value = 0.2998220186967728
$ws.Range("M17").Value = 0.4488921123198679
$ws.Range("B18").Value = 2.181880302188006
$ws.Range("C18").Value = 0.4122891199007768
$ws.Range("E18").Value = 0.06504758517637832
$ws.Range("F18").Value = 3.629700279093441
$ws.Range("G18").Value = 0.002571284423762977
$ws.Range("I18").Value = 2.398232131911982
$ws.Range("J18").Value = 0.1779937965088791
$ws.Range("L18").Value = 0.2996780893364388
$ws.Range("M18").Value = 0.4460067011231317
$ws.Range("B19").Value = 2.174461792563136
$ws.Range("C19").Value = 0.4091608508239233
$ws.Range("E19").Value = 0.0649138298131362
$ws.Range("F19").Value = 3.62573875278818
$ws.Range("G19").Value = 0.00257174444107132
$ws.Range("I19").Value = 2.395624637272704
$ws.Range("J19").Value = 0.1779849529156365
$ws.Range("L19").Value = 0.2996314920296399
$ws.Range("M19").Value = 0.4450343247261017
$ws.Range("B20").Value = 2.207938651767904
$ws.Range("C20").Value = 0.4232502218553691
$ws.Range("E20").Value = 0.06551506377695659
$ws.Range("F20").Value = 3.643671159041219
$ws.Range("G20").Value = 0.002569686668615245
$ws.Range("I20").Value = 2.40742397639842
$ws.Range("J20").Value = 0.178026922509094
$ws.Range("L20").Value = 0.2998496699737032
$ws.Range("M20").Value = 0.4494283099736833
$ws.Range("B21").Value = 2.322042928087058
$ws.Range("C21").Value = 0.4707996072121432
$ws.Range("E21").Value = 0.0675235724427381
$ws.Range("F21").Value = 3.705767479668054
$ws.Range("G21").Value = 0.002562988410794276
$ws.Range("I21").Value = 2.44821765721754
$ws.Range("J21").Value = 0.178206301683197
$ws.Range("L21").Value = 0.3007307397768244
$ws.Range("M21").Value = 0.4645097628089516
$ws.Range("B22").Value = 2.397833572877516
$ws.Range("C22").Value = 0.5020384381834333
$ws.Range("E22").Value = 0.0688280585989105
$ws.Range("F22").Value = 3.747734938054009
$ws.Range("G22").Value = 0.002558770575425991
$ws.Range("I22").Value = 2.475743174860384
$ws.Range("J22").Value = 0.1783525565595099
$ws.Range("L22").Value = 0.301416335390023
$ws.Range("M22").Value = 0.4746036596884409
$ws.Range("B23").Value = 2.357271031364462
$ws.Range("C23").Value = 0.4853504705288287
$ws.Range("E23").Value = 0.06813255196493984
$ws.Range("F23").Value = 3.725208959621853
$ws.Range("G23").Value = 0.002561007137661295
$ws.Range("I23").Value = 2.460972688225269
$ws.Range("J23").Value = 0.1782718148547744
$ws.Range("L23").Value = 0.3010404120363219
$ws.Range("M23").Value = 0.4691946374922793
$ws.Range("B24").Value = 2.206094154541745
$ws.Range("C24").Value = 0.4224757203529066
$ws.Range("E24").Value = 0.06548209141905481
$ws.Range("F24").Value = 3.642679473313649
$ws.Range("G24").Value = 0.002569798858569258
$ws.Range("I24").Value = 2.406771708180329
$ws.Range("J24").Value = 0.1780244744798871
$ws.Range("L24").Value = 0.2998371303902729
$ws.Range("M24").Value = 0.4491858161566853
$ws.Range("B25").Value = 2.048220344868128
$ws.Range("C25").Value = 0.3553187732978245
$ws.Range("E25").Value = 0.06258524294747048
$ws.Range("F25").Value = 3.559546699393536
$ws.Range("G25").Value = 0.002579973778448587
$ws.Range("I25").Value = 2.351968455578728
$ws.Range("J25").Value = 0.1778795231868671
$ws.Range("L25").Value = 0.2990140823312046
$ws.Range("M25").Value = 0.4286212772320681
